{"js": "// Append Steven Gsell's icebreaker entry to the end of the document:\n//   1) a blank paragraph\n//   2) a paragraph with the date + name\n//   3) a multi-run bio paragraph (with a couple of proofErr gram markers,\n//      mirroring the Word grammar-checker artifacts in the source XML)\n\n// Text runs for the bio paragraph, in order. Two of them (\"gain \" / \"that\"\n// and \"really just\") are wrapped by a grammar-check proofErr pair in the\n// canonical XML, recorded via the gramStart/gramEnd marker objects below.\nconst bioRuns = [\n  { text: \"Hello\" },\n  { text: \"! I\\u2019m Steven Gsell. I was born in Clear Water Florida but moved to Jacksonville as a wee lad (age 4)\" },\n  { text: \" and have lived here ever since. I hope to be able to travel around the US and find a place to call home with my wife\" },\n  { text: \". \" },\n  { text: \"I am currently in my last semester for my A.S in Computer Information Technology\" },\n  { text: \" and hope to further my education to a B.S. \" },\n  { text: \"I have \" },\n  { text: \"taught myself \" },\n  { text: \"a good bit with JavaScript\" },\n  { text: \", HTML, CSS, React, and MongoDB. I am really enjoying Java and hope to \" },\n  { gramStart: true },\n  { text: \"gain \" },\n  { text: \"that\" },\n  { gramEnd: true },\n  { text: \" much more knowledge for my toolbelt. I like\" },\n  { text: \" 3D printing, low voltage tinkering, web design/development\" },\n  { text: \", \" },\n  { text: \"gaming, movies, \" },\n  { text: \"whittling, crocheting (I \" },\n  { gramStart: true },\n  { text: \"really just\" },\n  { gramEnd: true },\n  { text: \" like to make stuff, physical, digital, or both!).\" },\n];\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nfunction runXml(text) {\n  // Mirror Word's own behavior: space-padded text keeps xml:space=\"preserve\".\n  const preserve = text !== text.trim() || text.length === 0;\n  const spaceAttr = preserve ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:r><w:t${spaceAttr}>${xmlEscape(text)}</w:t></w:r>`;\n}\n\nfunction paragraphOoxmlPackage(innerParagraphXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>${innerParagraphXml}</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n}\n\nconst body = context.document.body;\n\n// 1) Blank paragraph at the very end of the document. Route it through\n// insertOoxml too so it lands as a truly empty <w:p/> (insertParagraph(\"\")\n// otherwise leaves a stray empty <w:r/> inside it).\nconst blankPara = body.insertParagraph(\"\", Word.InsertLocation.end);\nblankPara.insertOoxml(paragraphOoxmlPackage(\"<w:p/>\"), Word.InsertLocation.replace);\n\n// 2) Date + name line.\nconst dateNamePara = body.insertParagraph(\"5/8/2023 Steven Gsell\", Word.InsertLocation.end);\n\n// 3) Bio paragraph, built with explicit runs (and proofErr markers) via\n// insertOoxml so formatting-identical runs don't get silently merged and\n// the gramStart/gramEnd pairs land exactly where the source XML has them.\nlet bioInner = \"\";\nfor (const token of bioRuns) {\n  if (token.gramStart) {\n    bioInner += '<w:proofErr w:type=\"gramStart\"/>';\n  } else if (token.gramEnd) {\n    bioInner += '<w:proofErr w:type=\"gramEnd\"/>';\n  } else {\n    bioInner += runXml(token.text);\n  }\n}\nconst bioParagraphXml = `<w:p>${bioInner}</w:p>`;\n\nconst bioPara = body.insertParagraph(\"\", Word.InsertLocation.end);\nbioPara.insertOoxml(paragraphOoxmlPackage(bioParagraphXml), Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Append Steven Gsell's icebreaker entry to the end of the document:\n#   1) a blank paragraph\n#   2) a paragraph with the date + name\n#   3) a multi-run bio paragraph (with a couple of proofErr gram markers,\n#      mirroring the Word grammar-checker artifacts in the source XML)\n#\n# Built as a single WordprocessingML fragment and inserted in one shot via\n# Range.InsertXML, so paragraph breaks land exactly where intended and the\n# proofErr markers stay attached to the right runs.\n\nfunction Esc([string]$s) {\n    $s = $s -replace '&', '&amp;'\n    $s = $s -replace '<', '&lt;'\n    $s = $s -replace '>', '&gt;'\n    return $s\n}\n\nfunction RunXml([string]$text) {\n    $preserveAttr = ''\n    if ($text -ne $text.Trim() -or $text.Length -eq 0) {\n        $preserveAttr = ' xml:space=\"preserve\"'\n    }\n    return '<w:r><w:t' + $preserveAttr + '>' + (Esc $text) + '</w:t></w:r>'\n}\n\n# Text runs for the bio paragraph, in order. \"gain \"/\"that\" and \"really just\"\n# are each wrapped by a grammar-check proofErr pair in the canonical XML.\n$bioRuns = @(\n    \"Hello\",\n    \"! I\u2019m Steven Gsell. I was born in Clear Water Florida but moved to Jacksonville as a wee lad (age 4)\",\n    \" and have lived here ever since. I hope to be able to travel around the US and find a place to call home with my wife\",\n    \". \",\n    \"I am currently in my last semester for my A.S in Computer Information Technology\",\n    \" and hope to further my education to a B.S. \",\n    \"I have \",\n    \"taught myself \",\n    \"a good bit with JavaScript\",\n    \", HTML, CSS, React, and MongoDB. I am really enjoying Java and hope to \"\n)\n$gramPair1 = @(\"gain \", \"that\")\n$bioRunsMid = @(\n    \" much more knowledge for my toolbelt. I like\",\n    \" 3D printing, low voltage tinkering, web design/development\",\n    \", \",\n    \"gaming, movies, \",\n    \"whittling, crocheting (I \"\n)\n$gramPair2 = @(\"really just\")\n$bioRunsEnd = @(\n    \" like to make stuff, physical, digital, or both!).\"\n)\n\n$bioInner = ''\nforeach ($t in $bioRuns) { $bioInner += RunXml $t }\n$bioInner += '<w:proofErr w:type=\"gramStart\"/>'\n$bioInner += RunXml $gramPair1[0]\n$bioInner += RunXml $gramPair1[1]\n$bioInner += '<w:proofErr w:type=\"gramEnd\"/>'\nforeach ($t in $bioRunsMid) { $bioInner += RunXml $t }\n$bioInner += '<w:proofErr w:type=\"gramStart\"/>'\n$bioInner += RunXml $gramPair2[0]\n$bioInner += '<w:proofErr w:type=\"gramEnd\"/>'\nforeach ($t in $bioRunsEnd) { $bioInner += RunXml $t }\n\n$blankParaXml = '<w:p/>'\n$dateNameParaXml = '<w:p>' + (RunXml \"5/8/2023 Steven Gsell\") + '</w:p>'\n$bioParaXml = '<w:p>' + $bioInner + '</w:p>'\n\n$bodyFragment = $blankParaXml + $dateNameParaXml + $bioParaXml\n\n$wordOpenXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $bodyFragment + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$d = $word.ActiveDocument\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertXML($wordOpenXml)\n"}
